# The deck's single writable theme part (ppt/theme/theme1.xml, "Integral")
# gets its 12 theme colors swapped for the "Office Theme" palette that
# previously only lived in ppt/theme/theme2.xml (used by the notes master).
#
# PowerPoint's ThemeColorScheme exposes exactly the 12 OOXML <a:clrScheme>
# slots in document order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# Item(1..12) on any slide/slide-range all reseat onto that same shared
# theme part, so we touch the scheme once via a SlideRange (no per-slide
# side effects) and write every slot to the "Office Theme" RGB values.

function ConvertTo-ComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" clrScheme (formerly theme2.xml).
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$p  = $ppt.ActivePresentation
$sr = $p.Slides.Range(1)
$tcs = $sr.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Item($i).RGB = ConvertTo-ComRGB $officeThemeColors[$i - 1]
}
